# Insert a new price-record row for "Feria Lagunitas de Puerto Montt - Betarraga"
# just above the current row 303, shifting the existing rows 303-317 down to
# 304-318 (dimension grows from A1:R317 to A1:R318).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 303..317 down by one to make room for the new record.
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new weekly record.
$ws.Range("A303").Value = 4
$ws.Range("B303").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C303").Value = "Los Lagos"
$ws.Range("D303").Value = 44753
$ws.Range("E303").Value = 10
$ws.Range("F303").Value = 100114014
$ws.Range("G303").Value = "Betarraga"
$ws.Range("H303").Value = "Sin especificar"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 500
$ws.Range("K303").Value = 1200
$ws.Range("L303").Value = 1300
$ws.Range("M303").Value = 1250
$ws.Range("N303").Value = "$/paquete 5 unidades"
$ws.Range("O303").Value = "Región del Maule"
$ws.Range("P303").Value = 250
$ws.Range("Q303").Value = 5
$ws.Range("R303").Value = "Hortaliza"
